# "admin buy product added" - replace the placeholder/test product names
# in column B (rows 2-16) with the new "Shop Pd - N" naming scheme.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($i = 1; $i -le 15; $i++) {
    $ws.Cells.Item($i + 1, 2).Value = "Shop Pd - $i"
}

# Move the active selection to F11 (matches the saved sheet view state).
$ws.Range("F11").Select()
